# Add the new "stocklvl" sheet (after "dcr"), populate it with the stockpile
# level time series, and update the active-sheet/selection state to match.

$wb = $excel.ActiveWorkbook

# --- Add sheet "stocklvl" as the last sheet (after "dcr") -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "stocklvl"
$ws.Activate()

# --- Header row -------------------------------------------------------------
$ws.Range("A1").Value = "Stf"
$ws.Range("B1").Value = "Value"

# --- Data rows (years 2024-2050 with their stockpile level values) --------
$years = @(2024, 2025, 2026, 2027, 2028, 2029, 2030, 2031, 2032, 2033, 2034, 2035, 2036, 2037, 2038, 2039, 2040, 2041, 2042, 2043, 2044, 2045, 2046, 2047, 2048, 2049, 2050)
$vals  = @(26229.707720000002, 26589.174669999997, 60431.69513, 95314.195130000007, 130196.6951, 165780.465, 179750.54879999999, 195610.7611, 196499.72339999999, 198853.9976, 202089.00579999998, 164077.8836, 126654.49769999999, 90701.344779999999, 57298.538379999998, 40938.322699999997, 30788.296410000003, 28368.054080000002, 31250.81408, 34011.350160000002, 35066.675019999995, 33453.543590000001, 25310.301890000002, 11877.08452, 0, 0, 0)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $vals[$i]
}

# New sheet's selection: single active cell H12.
$ws.Range("H12").Select()

# --- "dcr" sheet: selection becomes the whole used range, active cell B28 -
$dcr = $wb.Worksheets.Item("dcr")
$dcr.Activate()
$dcr.Range("A1:B28").Select()

# --- Re-activate the new sheet so it ends up the active/selected tab ------
$ws.Activate()
